$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.059.85'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '3.152.02'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '601.37'
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").Value = '142.53'
$ws.Range("E6").Value = '  -2.76%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.148.16'
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("E13").Value = '  -1.90%  '
$ws.Range("D14").Value = '34.98'
$ws.Range("E14").Value = '  -2.66%  '
$ws.Range("D15").Value = '3.669.67'
$ws.Range("E15").Value = '  -0.74%  '
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '64.031.80'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").Value = '3.146.07'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").Value = '486.76'
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").Value = '14.67'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  -1.39%  '
$ws.Range("D23").Value = '7.74'
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("D24").Value = '88.21'
$ws.Range("E24").Value = '  +4.63%  '
$ws.Range("D25").Value = '13.25'
$ws.Range("E25").Value = '  -4.04%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '2.77'
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("D28").Value = '8.26'
$ws.Range("E28").Value = '  -5.29%  '
$ws.Range("D29").Value = '6.99'
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("E30").Value = '  -2.51%  '
$ws.Range("D31").Value = '27.45'
$ws.Range("E31").Value = '  +2.90%  '
$ws.Range("E32").Value = '  -6.04%  '
$ws.Range("E33").Value = '  -0.23%  '
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E35").Value = '  -2.43%  '
$ws.Range("D36").Value = '6.07'
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '0.0₃0745'
$ws.Range("E38").Value = '  -5.77%  '
$ws.Range("D39").Value = '2.94'
$ws.Range("E39").Value = '  -8.11%  '
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("D41").Value = '433.77'
$ws.Range("E41").Value = '  -6.37%  '
$ws.Range("E42").Value = '  -0.32%  '
$ws.Range("D43").Value = '8.39'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("D44").Value = '2.917.18'
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("D45").Value = '0.261'
$ws.Range("E45").Value = '  -3.03%  '
$ws.Range("E46").Value = '  -5.90%  '
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").Value = '25.87'
$ws.Range("E49").Value = '  -3.33%  '
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").Value = '121.13'
$ws.Range("E51").Value = '  +0.33%  '